$d = $word.ActiveDocument

# Table 1
$t1 = $d.Tables.Item(1)
$t1.Cell(2,3).Range.Text = "2.033"
$t1.Cell(2,5).Range.Text = "4.748"
$t1.Cell(2,6).Range.Text = "0.077"
$t1.Cell(3,3).Range.Text = "4.536"
$t1.Cell(3,4).Range.Text = "0.007"
$t1.Cell(3,5).Range.Text = "4.045"
$t1.Cell(3,6).Range.Text = "0.0165"
$t1.Cell(4,5).Range.Text = "91.207"

# Table 2
$t2 = $d.Tables.Item(2)
$t2.Cell(2,3).Range.Text = "2.052"
$t2.Cell(2,5).Range.Text = "4.770"
$t2.Cell(2,6).Range.Text = "0.076"
$t2.Cell(3,3).Range.Text = "4.628"
$t2.Cell(3,4).Range.Text = "0.007"
$t2.Cell(3,5).Range.Text = "4.141"
$t2.Cell(3,6).Range.Text = "0.0155"
$t2.Cell(4,5).Range.Text = "91.088"

# Table 3
$t3 = $d.Tables.Item(3)
$t3.Cell(2,3).Range.Text = "23.578"
$t3.Cell(3,3).Range.Text = "0.245"
$t3.Cell(3,4).Range.Text = "0.621"

# Table 4
$t4 = $d.Tables.Item(4)
$t4.Cell(2,3).Range.Text = "1.907"
$t4.Cell(2,4).Range.Text = "0.008"
$t4.Cell(2,5).Range.Text = "4.631"
$t4.Cell(2,6).Range.Text = "0.0835"
$t4.Cell(3,3).Range.Text = "3.010"
$t4.Cell(3,5).Range.Text = "3.217"
$t4.Cell(3,6).Range.Text = "0.0415"
$t4.Rows.Item(4).Height = 28.55
$t4.Cell(4,5).Range.Text = "92.152"

# Table 5
$t5 = $d.Tables.Item(5)
$t5.Cell(2,3).Range.Text = "23.438"
$t5.Cell(3,3).Range.Text = "4.235"
$t5.Cell(3,4).Range.Text = "0.04*"
